$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("LiveData")

$ws.Range("B2").Value = 7193
$ws.Range("C3").Value = 172582
$ws.Range("C4").Value = 163373
$ws.Range("C8").Value = 66.01000000000001
